$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the LR-pair table (Ntn1-Mcam) to the recomputed 3-cluster x 3-cluster
# (ECs / FAPs / sCs) cross-tabulation following Dr Hou's advice.
# Rows 2-10 (A:T) are rewritten cell-by-cell with the new sending/target
# cluster labels and recalculated statistics; the table grows from 6 to 9
# data rows (dimension A1:T7 -> A1:T10).

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ntn1"
$ws.Cells.Item(2,3).Value = "Mcam"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.66666666666666663
$ws.Cells.Item(2,7).Value = 0.9305633333333333
$ws.Cells.Item(2,8).Value = 2.79169
$ws.Cells.Item(2,9).Value = 0.01768777137856805
$ws.Cells.Item(2,10).Value = 0.01768777137856806
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 84.077541666666676
$ws.Cells.Item(2,14).Value = 252.23262500000001
$ws.Cells.Item(2,15).Value = 0.51957148007956833
$ws.Cells.Item(2,16).Value = 0.51957148007956833
$ws.Cells.Item(2,17).Value = 78.239477431805554
$ws.Cells.Item(2,18).Value = 704.15529688625008
$ws.Cells.Item(2,19).Value = 0.0091900615544716294
$ws.Cells.Item(2,20).Value = 0.0091900615544716312
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ntn1"
$ws.Cells.Item(3,3).Value = "Mcam"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.66666666666666663
$ws.Cells.Item(3,7).Value = 0.9305633333333333
$ws.Cells.Item(3,8).Value = 2.79169
$ws.Cells.Item(3,9).Value = 0.01768777137856805
$ws.Cells.Item(3,10).Value = 0.01768777137856806
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.71712600000000004
$ws.Cells.Item(3,14).Value = 2.1513779999999998
$ws.Cells.Item(3,15).Value = 0.0044316021833837784
$ws.Cells.Item(3,16).Value = 0.0044316021833837784
$ws.Cells.Item(3,17).Value = 0.66733116098
$ws.Cells.Item(3,18).Value = 6.0059804488200008
$ws.Cells.Item(3,19).Value = 0.00007838516626045528
$ws.Cells.Item(3,20).Value = 0.0000783851662604553
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ntn1"
$ws.Cells.Item(4,3).Value = "Mcam"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.66666666666666663
$ws.Cells.Item(4,7).Value = 0.9305633333333333
$ws.Cells.Item(4,8).Value = 2.79169
$ws.Cells.Item(4,9).Value = 0.01768777137856805
$ws.Cells.Item(4,10).Value = 0.01768777137856806
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 77.026265333333342
$ws.Cells.Item(4,14).Value = 231.07879600000001
$ws.Cells.Item(4,15).Value = 0.47599691773704778
$ws.Cells.Item(4,16).Value = 0.47599691773704789
$ws.Cells.Item(4,17).Value = 71.677818222804447
$ws.Cells.Item(4,18).Value = 645.10036400524007
$ws.Cells.Item(4,19).Value = 0.0084193246578359669
$ws.Cells.Item(4,20).Value = 0.0084193246578359686
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ntn1"
$ws.Cells.Item(5,3).Value = "Mcam"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 44.154176999999997
$ws.Cells.Item(5,8).Value = 132.46253100000001
$ws.Cells.Item(5,9).Value = 0.83926473374711519
$ws.Cells.Item(5,10).Value = 0.8392647337471153
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 84.077541666666676
$ws.Cells.Item(5,14).Value = 252.23262500000001
$ws.Cells.Item(5,15).Value = 0.51957148007956833
$ws.Cells.Item(5,16).Value = 0.51957148007956833
$ws.Cells.Item(5,17).Value = 3712.3746564748758
$ws.Cells.Item(5,18).Value = 33411.371908273883
$ws.Cells.Item(5,19).Value = 0.4360580198915735
$ws.Cells.Item(5,20).Value = 0.43605801989157361
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ntn1"
$ws.Cells.Item(6,3).Value = "Mcam"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 44.154176999999997
$ws.Cells.Item(6,8).Value = 132.46253100000001
$ws.Cells.Item(6,9).Value = 0.83926473374711519
$ws.Cells.Item(6,10).Value = 0.8392647337471153
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.71712600000000004
$ws.Cells.Item(6,14).Value = 2.1513779999999998
$ws.Cells.Item(6,15).Value = 0.0044316021833837784
$ws.Cells.Item(6,16).Value = 0.0044316021833837784
$ws.Cells.Item(6,17).Value = 31.664108335302
$ws.Cells.Item(6,18).Value = 284.97697501771808
$ws.Cells.Item(6,19).Value = 0.0037192874265107199
$ws.Cells.Item(6,20).Value = 0.0037192874265107221
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ntn1"
$ws.Cells.Item(7,3).Value = "Mcam"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 44.154176999999997
$ws.Cells.Item(7,8).Value = 132.46253100000001
$ws.Cells.Item(7,9).Value = 0.83926473374711519
$ws.Cells.Item(7,10).Value = 0.8392647337471153
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 77.026265333333342
$ws.Cells.Item(7,14).Value = 231.07879600000001
$ws.Cells.Item(7,15).Value = 0.47599691773704778
$ws.Cells.Item(7,16).Value = 0.47599691773704789
$ws.Cells.Item(7,17).Value = 3401.0313531769648
$ws.Cells.Item(7,18).Value = 30609.282178592679
$ws.Cells.Item(7,19).Value = 0.39948742642903101
$ws.Cells.Item(7,20).Value = 0.39948742642903101
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Ntn1"
$ws.Cells.Item(8,3).Value = "Mcam"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 7.5258070000000004
$ws.Cells.Item(8,8).Value = 22.577421000000001
$ws.Cells.Item(8,9).Value = 0.14304749487431681
$ws.Cells.Item(8,10).Value = 0.14304749487431681
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 84.077541666666676
$ws.Cells.Item(8,14).Value = 252.23262500000001
$ws.Cells.Item(8,15).Value = 0.51957148007956833
$ws.Cells.Item(8,16).Value = 0.51957148007956833
$ws.Cells.Item(8,17).Value = 632.75135161779178
$ws.Cells.Item(8,18).Value = 5694.7621645601257
$ws.Cells.Item(8,19).Value = 0.074323398633523233
$ws.Cells.Item(8,20).Value = 0.074323398633523233
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Ntn1"
$ws.Cells.Item(9,3).Value = "Mcam"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 7.5258070000000004
$ws.Cells.Item(9,8).Value = 22.577421000000001
$ws.Cells.Item(9,9).Value = 0.14304749487431681
$ws.Cells.Item(9,10).Value = 0.14304749487431681
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.71712600000000004
$ws.Cells.Item(9,14).Value = 2.1513779999999998
$ws.Cells.Item(9,15).Value = 0.0044316021833837784
$ws.Cells.Item(9,16).Value = 0.0044316021833837784
$ws.Cells.Item(9,17).Value = 5.3969518706820008
$ws.Cells.Item(9,18).Value = 48.572566836138009
$ws.Cells.Item(9,19).Value = 0.00063392959061260199
$ws.Cells.Item(9,20).Value = 0.0006339295906126021
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Ntn1"
$ws.Cells.Item(10,3).Value = "Mcam"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 7.5258070000000004
$ws.Cells.Item(10,8).Value = 22.577421000000001
$ws.Cells.Item(10,9).Value = 0.14304749487431681
$ws.Cells.Item(10,10).Value = 0.14304749487431681
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 77.026265333333342
$ws.Cells.Item(10,14).Value = 231.07879600000001
$ws.Cells.Item(10,15).Value = 0.47599691773704778
$ws.Cells.Item(10,16).Value = 0.47599691773704789
$ws.Cells.Item(10,17).Value = 579.68480682945744
$ws.Cells.Item(10,18).Value = 5217.1632614651162
$ws.Cells.Item(10,19).Value = 0.068090166650180931
$ws.Cells.Item(10,20).Value = 0.068090166650180944
